$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 23) down to the five
# new rows (24-28) so the new cells pick up the same styles (date / time
# number formats) already used by the table, instead of creating new style
# entries.
$ws.Range("A23:D23").Copy()
$ws.Range("A24:D28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 24: 16/10/2023, 17:30 - 19:30
$ws.Range("A24").Value = 45215
$ws.Range("B24").Value = 0.72916666666666663
$ws.Range("C24").Value = 0.8125
$ws.Range("D24").Value = "All"
$ws.Range("E24").Value = "Agreed to complete task 1 by 20/10/2023 and meet up again for the findings"

# Row 25: 20/10/2023, 15:00 - 17:00
$ws.Range("A25").Value = 45219
$ws.Range("B25").Value = 0.625
$ws.Range("C25").Value = 0.70833333333333337
$ws.Range("D25").Value = "All"
$ws.Range("E25").Value = "Discuss findings of task 1 and set deadline for task 2 on 23/10/2023"

# Row 26: 23/10/2023, 19:00 - 21:00
$ws.Range("A26").Value = 45222
$ws.Range("B26").Value = 0.79166666666666663
$ws.Range("C26").Value = 0.875
$ws.Range("D26").Value = "All"
$ws.Range("E26").Value = "Discuss findings of task 2 and discuss on how to write up introduction and data description "

# Row 27: 25/10/2023, 20:00 - 21:00
$ws.Range("A27").Value = 45224
$ws.Range("B27").Value = 0.83333333333333337
$ws.Range("C27").Value = 0.875
$ws.Range("D27").Value = "All"
$ws.Range("E27").Value = "Discuss on minor errors on visual inference and finalise the summary part"

# Row 28: 26/10/2023, 15:45 - 0:00
$ws.Range("A28").Value = 45225
$ws.Range("B28").Value = 0.65625
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = "All"
$ws.Range("E28").Value = "Finalise the entire report and film the presentation with the other group"

# Update the active selection to match the new last-entered cell.
$ws.Range("C28").Select() | Out-Null
